# "big fix to smoothness parameter"
# Updates the cmax camera-right value and appends four new morphological
# smoothing parameters (trimedgeof, openradius, closewidth, closeheight)
# to the camparam sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Best-effort: widen the workbook tab-bar ratio like the source commit did.
# (No-op on hosts that don't expose window chrome via COM.)
try {
    $wb.ActiveWindow.TabRatio = 0.713
} catch {
}

# cmax right-camera value: 6500 -> 7000
$ws.Range("C6").Value = 7000

# New parameter rows appended below medfiltsize (row 18)
$ws.Range("A19").Value = "trimedgeof"
$ws.Range("B19").Value = 5
$ws.Range("C19").Value = 5

$ws.Range("A20").Value = "openradius"
$ws.Range("B20").Value = 4
$ws.Range("C20").Value = 4

$ws.Range("A21").Value = "closewidth"
$ws.Range("B21").Value = 5
$ws.Range("C21").Value = 5

$ws.Range("A22").Value = "closeheight"
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 1

# Match the author's final selection
[void]$ws.Range("C22").Select()
